$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Hydrogen): D3 corrected -> now blank (was 567.440869632662)
$ws.Range("D3").ClearContents()

# Row 4 (Methanol): C4 corrected to 0 (was 14630.53104371214)
$ws.Range("C4").Value = 0

# Row 5 (Ammonia): C5 corrected to 0 (was 59720.51281601335)
$ws.Range("C5").Value = 0

# Row 7: relabel "Other" -> "Biogas" and correct its D value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 690.30559251208

# Row 8 (new): re-add the "Other" row (previously row 7) with a freshly
# corrected value, matching the table's full A:D row layout/styling.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)   # xlPasteFormats - reuse the label style
$ws.Range("A8").Value = "Other"

$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)

$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)

$ws.Range("D8").Value = 480.2460036365646
